# Insert a new weekly price record for "Femacal de La Calera" (Mango) at
# row 111, pushing the existing rows 111:182 down to 112:183.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 111 (shifts rows 111:182 -> 112:183).
$ws.Rows.Item(111).Insert()

# Populate the new row 111 with the new record's data.
$ws.Cells.Item(111, 1).Value = 3
$ws.Cells.Item(111, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(111, 3).Value = "Coquimbo"
$ws.Cells.Item(111, 4).Value = 44438
$ws.Cells.Item(111, 5).Value = 5
$ws.Cells.Item(111, 6).Value = "Fruta"
$ws.Cells.Item(111, 7).Value = 100108
$ws.Cells.Item(111, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(111, 9).Value = 100108002
$ws.Cells.Item(111, 10).Value = "Mango"
$ws.Cells.Item(111, 11).Value = "Sin especificar"
$ws.Cells.Item(111, 12).Value = "Primera"
$ws.Cells.Item(111, 13).Value = 456
$ws.Cells.Item(111, 14).Value = 9000
$ws.Cells.Item(111, 15).Value = 9000
$ws.Cells.Item(111, 16).Value = 9000
$ws.Cells.Item(111, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(111, 18).Value = "Brasil"
$ws.Cells.Item(111, 19).Value = 2250
$ws.Cells.Item(111, 20).Value = 4
